# Finished optimizing import from .xls
#
# - Appends short abbreviations to three oblast (region) names that are
#   repeated down column A for every hospital record belonging to that
#   region.
# - Updates the sheet's view state (scroll position / selection) to
#   reflect where the author was last working in the sheet.
# - Widens column A so the longer region names fit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update region names -------------------------------------------------
# Every row that belongs to a given oblast repeats the same region string
# in column A, so the whole contiguous block is rewritten at once (this is
# equivalent to editing the single shared string all of those cells point
# to).

# Восточно-Казахстанская область -> Восточно-Казахстанская область (ВКО)
$ws.Range("A987:A1211").Value = "Восточно-Казахстанская область (ВКО)"

# Западно-Казахстанская область -> Западно-Казахстанская область (ЗКО)
$ws.Range("A1836:A1987").Value = "Западно-Казахстанская область (ЗКО)"

# Северо-Казахстанская область -> Северо-Казахстанская область (СКО)
$ws.Range("A2661:A2734").Value = "Северо-Казахстанская область (СКО)"

# --- 2. Restore the author's scroll position / selection -------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 43
$win.ScrollColumn = 1

# The author had the A987:A1211 block (the region just renamed) selected.
$ws.Range("A987:A1211").Select()

# --- 3. Widen column A to fit the longer region names -----------------------
$ws.Columns.Item(1).ColumnWidth = 44.3
